# The deck ships two embedded themes:
#   ppt/theme/theme1.xml -> bound to the slide master ("Integral" palette)
#   ppt/theme/theme2.xml -> bound to the notes master  ("Office Theme" palette)
# The authored change swaps the two themes' colour palettes (the slide master
# picks up the stock "Office Theme" colours, the notes master gets "Integral").
# PowerPoint's object model only exposes the slide-master-bound theme colours
# for programmatic editing (Slide/SlideRange/Master.Theme.ThemeColorScheme);
# recolour it here to the "Office Theme" palette taken from theme2.xml.

function Hex-ToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Office Theme colour scheme, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1..6, hlink, folHlink
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $themeColors.Colors($i).RGB = Hex-ToVbaRgb $officeThemeHex[$i - 1]
}
